# Apply updated cryptos list values (prices + 1h volume %) scraped on
# Wed Jan 17 04:51:49 UTC 2024, including the ordi/BitcoinSV row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.876.31"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.569.10"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.85%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "2.958.53"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.20%  "
$ws.Range("D16").Value = "2.558.05"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "42.859.00"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("E38").Value = "  +10.98%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.15%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0303"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "2.006.72"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "2.810.69"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.197"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
